$d = $word.ActiveDocument

# --- "Caso de Uso" header row -------------------------------------------
# "UC00" + "3" were split across two runs; normalize to a single run with
# the same visible text "UC003" (no wording change, just a run merge).
$d.Content.Find.Execute(
    "UC003", $true, $false, $false, $false, $false,
    $true, 1, $false, "UC003", 2)

# --- "Atributos" row ------------------------------------------------------
# "Título do conteúdo, descrição do cont" + "eú" + "do, data do conteúdo"
# merged into one run with identical text.
$d.Content.Find.Execute(
    "Título do conteúdo, descrição do conteúdo, data do conteúdo", $true, $false, $false, $false, $false,
    $true, 1, $false, "Título do conteúdo, descrição do conteúdo, data do conteúdo", 2)

# --- "Fluxo Principal" row --------------------------------------------------
# Step 1: the user now goes to the content search screen instead of home.
$d.Content.Find.Execute(
    "usuário acessa sua home do aplicativo, após o login.", $true, $false, $false, $false, $false,
    $true, 1, $false, "usuário acessa sua tela de busca por conteúdos.", 2)

# Step 2: the app now "exibe" (shows) instead of "envia" (sends) suggestions.
$d.Content.Find.Execute(
    "aplicativo envia sugestão de conteúdos por meio de notificações push.", $true, $false, $false, $false, $false,
    $true, 1, $false, "aplicativo exibe sugestão de conteúdos por meio de notificações push.", 2)

# --- "Fluxo Alternativo" row ------------------------------------------------
# "buscar por sugestões de conte" + "ú" + "do clicando no botão sugestões."
# merged into one run with identical text.
$d.Content.Find.Execute(
    "buscar por sugestões de conteúdo clicando no botão sugestões.", $true, $false, $false, $false, $false,
    $true, 1, $false, "buscar por sugestões de conteúdo clicando no botão sugestões.", 2)

# "o aplicativo " + "solicita para o usuário ..." merged into one run with
# identical text.
$d.Content.Find.Execute(
    "o aplicativo solicita para o usuário qual o tipo de conteúdo que ele tem interesse e o aplicativo retorna uma lista com as informações sobre os conteúdos solicitados.", $true, $false, $false, $false, $false,
    $true, 1, $false, "o aplicativo solicita para o usuário qual o tipo de conteúdo que ele tem interesse e o aplicativo retorna uma lista com as informações sobre os conteúdos solicitados.", 2)
